$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

    $ws.Range("C16").Value = "1047462094"
    $ws.Range("D16").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E16").Value = "2205"
    $ws.Range("F16").Value = 40000
    $ws.Range("G16").Value = 1000000
    $ws.Range("C17").Value = "1047462094"
    $ws.Range("D17").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E17").Value = "2206"
    $ws.Range("F17").Value = 40000
    $ws.Range("G17").Value = 1000000
    $ws.Range("C18").Value = "1047462094"
    $ws.Range("D18").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E18").Value = "2207"
    $ws.Range("F18").Value = 40000
    $ws.Range("G18").Value = 1000000
    $ws.Range("C19").Value = "1047462094"
    $ws.Range("D19").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E19").Value = "2208"
    $ws.Range("F19").Value = 40000
    $ws.Range("G19").Value = 1000000
    $ws.Range("C20").Value = "1047462094"
    $ws.Range("D20").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E20").Value = "2209"
    $ws.Range("F20").Value = 40000
    $ws.Range("G20").Value = 1000000
    $ws.Range("C21").Value = "1047462094"
    $ws.Range("D21").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E21").Value = "2210"
    $ws.Range("F21").Value = 40000
    $ws.Range("G21").Value = 1000000
    $ws.Range("C22").Value = "1043588724"
    $ws.Range("D22").Value = "RONAL ALBERTO BELTRAN RAMIREZ"
    $ws.Range("E22").Value = "2210"
    $ws.Range("F22").Value = 32000
    $ws.Range("G22").Value = 1423500
    $ws.Range("C23").Value = "1047462094"
    $ws.Range("D23").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E23").Value = "2211"
    $ws.Range("F23").Value = 40000
    $ws.Range("G23").Value = 1000000
    $ws.Range("C24").Value = "1047462094"
    $ws.Range("D24").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E24").Value = "2212"
    $ws.Range("F24").Value = 40000
    $ws.Range("G24").Value = 1000000
    $ws.Range("C25").Value = "1047462094"
    $ws.Range("D25").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E25").Value = "2301"
    $ws.Range("F25").Value = 40000
    $ws.Range("G25").Value = 1000000
    $ws.Range("C26").Value = "3828322"
    $ws.Range("D26").Value = "WILSON JOSE RINCON ARELLANO"
    $ws.Range("E26").Value = "2302"
    $ws.Range("F26").Value = 17333
    $ws.Range("G26").Value = 1300000
    $ws.Range("C27").Value = "1047462094"
    $ws.Range("D27").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E27").Value = "2302"
    $ws.Range("F27").Value = 40000
    $ws.Range("G27").Value = 1000000
    $ws.Range("C28").Value = "1047462094"
    $ws.Range("D28").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E28").Value = "2303"
    $ws.Range("F28").Value = 40000
    $ws.Range("G28").Value = 1000000
    $ws.Range("C29").Value = "1047462094"
    $ws.Range("D29").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E29").Value = "2304"
    $ws.Range("F29").Value = 40000
    $ws.Range("G29").Value = 1000000
    $ws.Range("C30").Value = "1047462094"
    $ws.Range("D30").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E30").Value = "2305"
    $ws.Range("F30").Value = 40000
    $ws.Range("G30").Value = 1000000
    $ws.Range("C31").Value = "1047462094"
    $ws.Range("D31").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E31").Value = "2306"
    $ws.Range("F31").Value = 40000
    $ws.Range("G31").Value = 1000000
    $ws.Range("C32").Value = "1047462094"
    $ws.Range("D32").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E32").Value = "2307"
    $ws.Range("F32").Value = 40000
    $ws.Range("G32").Value = 1000000
    $ws.Range("C33").Value = "1047462094"
    $ws.Range("D33").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E33").Value = "2308"
    $ws.Range("F33").Value = 40000
    $ws.Range("G33").Value = 1000000
    $ws.Range("C34").Value = "1047462094"
    $ws.Range("D34").Value = "GUSTAVO ENRIQUE PAJARO AREVALO"
    $ws.Range("E34").Value = "2309"
    $ws.Range("F34").Value = 5333
    $ws.Range("G34").Value = 1000000
